$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '70.347.92'
Set-TextValue 'E2' '  +5.02%  '
Set-TextValue 'D3' '3.606.67'
Set-TextValue 'E3' '  +4.75%  '
Set-TextValue 'E4' '  +0.14%  '
Set-TextValue 'D5' '589.08'
Set-TextValue 'E5' '  +3.25%  '
Set-TextValue 'D6' '190.25'
Set-TextValue 'E6' '  +3.15%  '
Set-TextValue 'E7' '  +1.78%  '
Set-TextValue 'D8' '3.594.64'
Set-TextValue 'E8' '  +4.60%  '
Set-TextValue 'E9' '  +0.04%  '
Set-TextValue 'E10' '  +0.33%  '
Set-TextValue 'E11' '  +2.61%  '
Set-TextValue 'D12' '58.04'
Set-TextValue 'E12' '  +4.87%  '
Set-TextValue 'E13' '  +3.27%  '
Set-TextValue 'D14' '9.78'
Set-TextValue 'E14' '  +4.29%  '
Set-TextValue 'D15' '4.178.23'
Set-TextValue 'E15' '  +4.93%  '
Set-TextValue 'D16' '3.608.08'
Set-TextValue 'E16' '  +5.38%  '
Set-TextValue 'D17' '19.37'
Set-TextValue 'E17' '  +4.40%  '
Set-TextValue 'D18' '70.246.86'
Set-TextValue 'E18' '  +5.36%  '
Set-TextValue 'D19' '12.47'
Set-TextValue 'E19' '  +3.59%  '
Set-TextValue 'E20' '  +0.26%  '
Set-TextValue 'E21' '  +4.12%  '
Set-TextValue 'D22' '494.18'
Set-TextValue 'E22' '  +5.07%  '
Set-TextValue 'D23' '17.27'
Set-TextValue 'E23' '  +15.97%  '
Set-TextValue 'D24' '5.37'
Set-TextValue 'E24' '  +7.47%  '
Set-TextValue 'D25' '4.46'
Set-TextValue 'E25' '  +6.47%  '
Set-TextValue 'D26' '90.78'
Set-TextValue 'E26' '  +1.26%  '
Set-TextValue 'E27' '  +5.05%  '
Set-TextValue 'D28' '11.08'
Set-TextValue 'E28' '  +1.11%  '
Set-TextValue 'D29' '9.46'
Set-TextValue 'E29' '  +5.88%  '
Set-TextValue 'D30' '32.34'
Set-TextValue 'E30' '  +2.58%  '
Set-TextValue 'D31' '7.57'
Set-TextValue 'E31' '  +8.45%  '
Set-TextValue 'E32' '  +5.21%  '
Set-TextValue 'D33' '618.45'
Set-TextValue 'E33' '  +6.15%  '
Set-TextValue 'D34' '0.117'
Set-TextValue 'E34' '  +6.67%  '
Set-TextValue 'D35' '65.25'
Set-TextValue 'E35' '  +3.81%  '
Set-TextValue 'D36' '0.0₃0817'
Set-TextValue 'E36' '  +6.25%  '
Set-TextValue 'E37' '  +3.59%  '
Set-TextValue 'D38' '38.09'
Set-TextValue 'E38' '  +4.01%  '
Set-TextValue 'E39' '  +0.14%  '
Set-TextValue 'E40' '  -1.15%  '
Set-TextValue 'D41' '3.62'
Set-TextValue 'E41' '  -0.98%  '
Set-TextValue 'D42' '3.303.19'
Set-TextValue 'E42' '  +5.49%  '
Set-TextValue 'E43' '  +5.98%  '
Set-TextValue 'D44' '0.0445'
Set-TextValue 'E44' '  +4.69%  '
Set-TextValue 'E45' '  +2.26%  '
Set-TextValue 'E46' '  +0.68%  '
Set-TextValue 'E47' '  +2.14%  '
Set-TextValue 'E48' '  +5.51%  '
Set-TextValue 'D49' '2.72'
Set-TextValue 'E49' '  -2.86%  '
Set-TextValue 'D50' '3.31'
Set-TextValue 'E50' '  +5.61%  '
Set-TextValue 'B51' 'Monero'
Set-TextValue 'C51' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D51' '143.09'
Set-TextValue 'E51' '  +1.02%  '
